$d = $word.ActiveDocument

function Get-BookmarkByName($doc, $name) {
    $count = $doc.Bookmarks.Count
    for ($i = 1; $i -le $count; $i++) {
        $b = $doc.Bookmarks.Item($i)
        if ($b.Name -eq $name) {
            return $b
        }
    }
    return $null
}

function Rename-HeadingBookmark($doc, $oldName, $newName, $newText, $pStyle) {
    $bm = Get-BookmarkByName $doc $oldName
    if ($bm -eq $null) {
        Write-Output "WARNING: bookmark '$oldName' not found"
        return
    }

    $paraStart = $bm.Start
    $paraEndExclusive = $bm.End + 1   # include the paragraph mark

    # Build a package-level OpenXML fragment: an empty "filler" paragraph
    # (to absorb the boundary with the preceding paragraph, so the old
    # bookmark - which collapses to zero width when its text is replaced -
    # does not linger) followed by the real replacement heading paragraph.
    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p/>' +
        '<w:p><w:pPr><w:pStyle w:val="' + $pStyle + '"/></w:pPr>' +
        '<w:bookmarkStart w:id="0" w:name="' + $newName + '"/>' +
        '<w:r><w:t xml:space="preserve">' + $newText + '</w:t></w:r>' +
        '<w:bookmarkEnd w:id="0"/></w:p>' +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    # Replace from the end of the previous paragraph (its paragraph mark)
    # through the end of this paragraph (its own paragraph mark).
    $rng = $doc.Range($paraStart - 1, $paraEndExclusive)
    $rng.InsertXML($xmlFrag)

    # The filler paragraph's own mark now sits right at $paraStart; delete
    # it to merge the (empty) filler back into the new heading paragraph,
    # which drops the leftover zero-width old bookmark without disturbing
    # the previous paragraph's own formatting/paragraph mark.
    $mergeRng = $doc.Range($paraStart, $paraStart + 1)
    $mergeRng.Delete()
}

Rename-HeadingBookmark $d "übung-1-1" "übung-1-plattform" "Übung 1: Plattform" "Heading3"
Rename-HeadingBookmark $d "übung-2-1" "übung-2-podcastbeschreibung" "Übung 2: Podcastbeschreibung" "Heading3"
Rename-HeadingBookmark $d "übung-3-1" "übung-3-coverbild" "Übung 3: Coverbild" "Heading3"

Write-Output "done"
